$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.64162015914917
$ws.Range("B1").Value = 3.891210317611694
$ws.Range("C1").Value = 2.690644025802612
$ws.Range("D1").Value = 0.8882414102554321
$ws.Range("E1").Value = 0.8626736402511597
